$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.024.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.019.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  -4.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.319.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("E13").Value = "  -4.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.738"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.000.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.955.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("E25").Value = "  -6.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("E27").Value = "  -5.54%  "
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -5.22%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.457.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "95.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0907"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.31%  "
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.208.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.45%  "
